# Generate Report for Handoff
# b.md has been handed off again: a new handoff xliff was generated for
# both the zh-cn and de-de locales, so the Overview sheet and each locale
# sheet need to reflect the refreshed status/timestamps/handoff file and
# the new "stale handback" error detail.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9378f215a8a0e13395494d3cf29b1eed34158dce/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/91c4c9f6a76774d395343a12a991bfb37e637515/e2e/b.md."

# --- Overview sheet: row 3 is b.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-02 20:44:16"

# --- zh-cn sheet: row 3 is b.md ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-02 20:44:11"
$zhcn.Range("P3").Value = $errorDetail

# --- de-de sheet: row 3 is b.md ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-02 20:44:16"
$dede.Range("P3").Value = $errorDetail
